$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"0.6868273333333333"
$ws.Range("H2").Value = [double]"2.060482"
$ws.Range("I2").Value = [double]"0.01130642661970366"
$ws.Range("J2").Value = [double]"0.01130642661970366"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"0.3465496666666667"
$ws.Range("N2").Value = [double]"1.039649"
$ws.Range("O2").Value = [double]"0.008996151488293185"
$ws.Range("P2").Value = [double]"0.008996151488293185"
$ws.Range("Q2").Value = [double]"0.2380197834242223"
$ws.Range("R2").Value = [double]"2.142178050818"
$ws.Range("S2").Value = [double]"0.0001017143266621248"
$ws.Range("T2").Value = [double]"0.0001017143266621248"

$ws.Range("G3").Value = [double]"0.6868273333333333"
$ws.Range("H3").Value = [double]"2.060482"
$ws.Range("I3").Value = [double]"0.01130642661970366"
$ws.Range("J3").Value = [double]"0.01130642661970366"
$ws.Range("N3").Value = [double]"66.23320799999999"
$ws.Range("O3").Value = [double]"0.5731203249593199"
$ws.Range("P3").Value = [double]"0.5731203249593199"
$ws.Range("Q3").Value = [double]"15.16359254291733"
$ws.Range("R3").Value = [double]"136.472332886256"
$ws.Range("S3").Value = [double]"0.006479942898413266"
$ws.Range("T3").Value = [double]"0.006479942898413265"

$ws.Range("G4").Value = [double]"0.6868273333333333"
$ws.Range("H4").Value = [double]"2.060482"
$ws.Range("I4").Value = [double]"0.01130642661970366"
$ws.Range("J4").Value = [double]"0.01130642661970366"
$ws.Range("M4").Value = [double]"16.08941833333333"
$ws.Range("N4").Value = [double]"48.268255"
$ws.Range("O4").Value = [double]"0.4176683996767803"
$ws.Range("P4").Value = [double]"0.4176683996767803"
$ws.Range("Q4").Value = [double]"11.05065228876778"
$ws.Range("R4").Value = [double]"99.45587059891001"
$ws.Range("S4").Value = [double]"0.004722337112314576"
$ws.Range("T4").Value = [double]"0.004722337112314576"

$ws.Range("G5").Value = [double]"0.6868273333333333"
$ws.Range("H5").Value = [double]"2.060482"
$ws.Range("I5").Value = [double]"0.01130642661970366"
$ws.Range("J5").Value = [double]"0.01130642661970366"
$ws.Range("M5").Value = [double]"0.008287000000000001"
$ws.Range("N5").Value = [double]"0.024861"
$ws.Range("O5").Value = [double]"0.0002151238756065334"
$ws.Range("P5").Value = [double]"0.0002151238756065334"
$ws.Range("Q5").Value = [double]"0.005691738111333334"
$ws.Range("R5").Value = [double]"0.051225643002"
$ws.Range("S5").Value = [double]"2.432282313691528E-06"
$ws.Range("T5").Value = [double]"2.432282313691528E-06"

$ws.Range("G6").Value = [double]"53.540432"
$ws.Range("I6").Value = [double]"0.8813728519762372"
$ws.Range("J6").Value = [double]"0.881372851976237"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"0.3465496666666667"
$ws.Range("N6").Value = [double]"1.039649"
$ws.Range("O6").Value = [double]"0.008996151488293185"
$ws.Range("P6").Value = [double]"0.008996151488293185"
$ws.Range("Q6").Value = [double]"18.55441886278934"
$ws.Range("R6").Value = [double]"166.989769765104"
$ws.Range("S6").Value = [double]"0.007928963694047235"
$ws.Range("T6").Value = [double]"0.007928963694047234"

$ws.Range("G7").Value = [double]"53.540432"
$ws.Range("I7").Value = [double]"0.8813728519762372"
$ws.Range("J7").Value = [double]"0.881372851976237"
$ws.Range("N7").Value = [double]"66.23320799999999"
$ws.Range("O7").Value = [double]"0.5731203249593199"
$ws.Range("P7").Value = [double]"0.5731203249593199"
$ws.Range("S7").Value = [double]"0.5051326953349435"
$ws.Range("T7").Value = [double]"0.5051326953349435"

$ws.Range("G8").Value = [double]"53.540432"
$ws.Range("I8").Value = [double]"0.8813728519762372"
$ws.Range("J8").Value = [double]"0.881372851976237"
$ws.Range("M8").Value = [double]"16.08941833333333"
$ws.Range("N8").Value = [double]"48.268255"
$ws.Range("O8").Value = [double]"0.4176683996767803"
$ws.Range("P8").Value = [double]"0.4176683996767803"
$ws.Range("Q8").Value = [double]"861.4344081953868"
$ws.Range("R8").Value = [double]"7752.90967375848"
$ws.Range("S8").Value = [double]"0.3681215886034747"
$ws.Range("T8").Value = [double]"0.3681215886034747"

$ws.Range("G9").Value = [double]"53.540432"
$ws.Range("I9").Value = [double]"0.8813728519762372"
$ws.Range("J9").Value = [double]"0.881372851976237"
$ws.Range("M9").Value = [double]"0.008287000000000001"
$ws.Range("N9").Value = [double]"0.024861"
$ws.Range("O9").Value = [double]"0.0002151238756065334"
$ws.Range("P9").Value = [double]"0.0002151238756065334"
$ws.Range("Q9").Value = [double]"0.4436895599840001"
$ws.Range("R9").Value = [double]"3.993206039856"
$ws.Range("S9").Value = [double]"0.0001896043437715116"
$ws.Range("T9").Value = [double]"0.0001896043437715116"

$ws.Range("G10").Value = [double]"6.476716"
$ws.Range("H10").Value = [double]"19.430148"
$ws.Range("I10").Value = [double]"0.1066185206043934"
$ws.Range("J10").Value = [double]"0.1066185206043934"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"0.3465496666666667"
$ws.Range("N10").Value = [double]"1.039649"
$ws.Range("O10").Value = [double]"0.008996151488293185"
$ws.Range("P10").Value = [double]"0.008996151488293185"
$ws.Range("Q10").Value = [double]"2.244503770894667"
$ws.Range("R10").Value = [double]"20.200533938052"
$ws.Range("S10").Value = [double]"0.0009591563628148316"
$ws.Range("T10").Value = [double]"0.0009591563628148316"

$ws.Range("G11").Value = [double]"6.476716"
$ws.Range("H11").Value = [double]"19.430148"
$ws.Range("I11").Value = [double]"0.1066185206043934"
$ws.Range("J11").Value = [double]"0.1066185206043934"
$ws.Range("N11").Value = [double]"66.23320799999999"
$ws.Range("O11").Value = [double]"0.5731203249593199"
$ws.Range("P11").Value = [double]"0.5731203249593199"
$ws.Range("Q11").Value = [double]"142.991225994976"
$ws.Range("R11").Value = [double]"1286.921033954784"
$ws.Range("S11").Value = [double]"0.06110524117547191"
$ws.Range("T11").Value = [double]"0.06110524117547191"

$ws.Range("G12").Value = [double]"6.476716"
$ws.Range("H12").Value = [double]"19.430148"
$ws.Range("I12").Value = [double]"0.1066185206043934"
$ws.Range("J12").Value = [double]"0.1066185206043934"
$ws.Range("M12").Value = [double]"16.08941833333333"
$ws.Range("N12").Value = [double]"48.268255"
$ws.Range("O12").Value = [double]"0.4176683996767803"
$ws.Range("P12").Value = [double]"0.4176683996767803"
$ws.Range("Q12").Value = [double]"104.2065931501933"
$ws.Range("R12").Value = [double]"937.85933835174"
$ws.Range("S12").Value = [double]"0.04453118687674283"
$ws.Range("T12").Value = [double]"0.04453118687674284"

$ws.Range("G13").Value = [double]"6.476716"
$ws.Range("H13").Value = [double]"19.430148"
$ws.Range("I13").Value = [double]"0.1066185206043934"
$ws.Range("J13").Value = [double]"0.1066185206043934"
$ws.Range("M13").Value = [double]"0.008287000000000001"
$ws.Range("N13").Value = [double]"0.024861"
$ws.Range("O13").Value = [double]"0.0002151238756065334"
$ws.Range("P13").Value = [double]"0.0002151238756065334"
$ws.Range("Q13").Value = [double]"0.05367254549200001"
$ws.Range("R13").Value = [double]"0.483052909428"
$ws.Range("S13").Value = [double]"2.293618936385216E-05"
$ws.Range("T13").Value = [double]"2.293618936385216E-05"

$ws.Range("E14").Value = [double]"3"
$ws.Range("F14").Value = [double]"1"
$ws.Range("G14").Value = [double]"0.04265633333333333"
$ws.Range("H14").Value = [double]"0.127969"
$ws.Range("I14").Value = [double]"0.0007022007996657373"
$ws.Range("J14").Value = [double]"0.0007022007996657372"
$ws.Range("K14").Value = [double]"3"
$ws.Range("L14").Value = [double]"1"
$ws.Range("M14").Value = [double]"0.3465496666666667"
$ws.Range("N14").Value = [double]"1.039649"
$ws.Range("O14").Value = [double]"0.008996151488293185"
$ws.Range("P14").Value = [double]"0.008996151488293185"
$ws.Range("Q14").Value = [double]"0.01478253809788889"
$ws.Range("R14").Value = [double]"0.133042842881"
$ws.Range("S14").Value = [double]"6.317104768993587E-06"
$ws.Range("T14").Value = [double]"6.317104768993586E-06"

$ws.Range("E15").Value = [double]"3"
$ws.Range("F15").Value = [double]"1"
$ws.Range("G15").Value = [double]"0.04265633333333333"
$ws.Range("H15").Value = [double]"0.127969"
$ws.Range("I15").Value = [double]"0.0007022007996657373"
$ws.Range("J15").Value = [double]"0.0007022007996657372"
$ws.Range("N15").Value = [double]"66.23320799999999"
$ws.Range("O15").Value = [double]"0.5731203249593199"
$ws.Range("P15").Value = [double]"0.5731203249593199"
$ws.Range("Q15").Value = [double]"0.9417552660613332"
$ws.Range("R15").Value = [double]"8.475797394551998"
$ws.Range("S15").Value = [double]"0.0004024455504911216"
$ws.Range("T15").Value = [double]"0.0004024455504911216"

$ws.Range("E16").Value = [double]"3"
$ws.Range("F16").Value = [double]"1"
$ws.Range("G16").Value = [double]"0.04265633333333333"
$ws.Range("H16").Value = [double]"0.127969"
$ws.Range("I16").Value = [double]"0.0007022007996657373"
$ws.Range("J16").Value = [double]"0.0007022007996657372"
$ws.Range("M16").Value = [double]"16.08941833333333"
$ws.Range("N16").Value = [double]"48.268255"
$ws.Range("O16").Value = [double]"0.4176683996767803"
$ws.Range("P16").Value = [double]"0.4176683996767803"
$ws.Range("Q16").Value = [double]"0.6863155915661111"
$ws.Range("R16").Value = [double]"6.176840324095"
$ws.Range("S16").Value = [double]"0.0002932870842481439"
$ws.Range("T16").Value = [double]"0.0002932870842481439"

$ws.Range("E17").Value = [double]"3"
$ws.Range("F17").Value = [double]"1"
$ws.Range("G17").Value = [double]"0.04265633333333333"
$ws.Range("H17").Value = [double]"0.127969"
$ws.Range("I17").Value = [double]"0.0007022007996657373"
$ws.Range("J17").Value = [double]"0.0007022007996657372"
$ws.Range("M17").Value = [double]"0.008287000000000001"
$ws.Range("N17").Value = [double]"0.024861"
$ws.Range("O17").Value = [double]"0.0002151238756065334"
$ws.Range("P17").Value = [double]"0.0002151238756065334"
$ws.Range("Q17").Value = [double]"0.0003534930343333333"
$ws.Range("R17").Value = [double]"0.003181437309"
$ws.Range("S17").Value = [double]"1.510601574781004E-07"
$ws.Range("T17").Value = [double]"1.510601574781004E-07"
